# DC-Colos.xlsx update
# Mirrors an upstream "generated data" refresh where the colo list was
# resorted. Net effect observed in the diff:
#   - Row 103 (CGB / Cuiaba, Brazil) is removed from its spot; rows 104-147
#     shift up by one row, and the old row-103 data is appended at row 147.
#   - Row 172 (RUN / Saint-Denis, Reunion) is removed from its spot; rows
#     173-179 shift up by one row, and the old row-172 data is appended at
#     row 179.
#   - Row 310 (FSD / Sioux Falls) gets its "name"/"country" text updated to
#     the newer naming convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Move-RowToEnd {
    param($sheet, $StartRow, $EndRow, $Columns)

    # Remember the row that is about to be displaced.
    $saved = @()
    for ($c = 1; $c -le $Columns; $c++) {
        $saved += $sheet.Cells.Item($StartRow, $c).Value()
    }

    # Shift every following row up by one.
    for ($r = $StartRow + 1; $r -le $EndRow; $r++) {
        for ($c = 1; $c -le $Columns; $c++) {
            $val = $sheet.Cells.Item($r, $c).Value()
            $sheet.Cells.Item($r - 1, $c).Value = $val
        }
    }

    # Drop the displaced row's data into the now-empty slot at the end.
    for ($c = 1; $c -le $Columns; $c++) {
        $sheet.Cells.Item($EndRow, $c).Value = $saved[$c - 1]
    }
}

# colo / name / region / city / country / cca2 / lat / lon => 8 columns
Move-RowToEnd $ws 103 147 8
Move-RowToEnd $ws 172 179 8

# Sioux Falls naming convention update.
$ws.Range("B310").Value = "Sioux Falls, SD, United States"
$ws.Range("E310").Value = "United States"
